{"js": "const pairs = [\n  [\"33\u00d765=2145\", \"83\u00d794=7802\"],\n  [\"86\u00d725=2150\", \"20\u00d780=1600\"],\n  [\"92\u00d787=8004\", \"65\u00d780=5200\"],\n  [\"87\u00d793=8091\", \"63\u00d741=2583\"],\n  [\"97\u00d733=3201\", \"34\u00d798=3332\"],\n  [\"17\u00d794=1598\", \"84\u00d772=6048\"],\n  [\"55\u00d794=5170\", \"12\u00d752=624\"],\n  [\"59\u00d717=1003\", \"45\u00d784=3780\"],\n  [\"57\u00d749=2793\", \"98\u00d790=8820\"],\n  [\"70\u00d742=2940\", \"29\u00d748=1392\"],\n  [\"73\u00d753=3869\", \"36\u00d760=2160\"],\n  [\"92\u00d762=5704\", \"97\u00d734=3298\"],\n  [\"32\u00d713=416\", \"25\u00d782=2050\"],\n  [\"11\u00d770=770\", \"65\u00d758=3770\"],\n  [\"73\u00d765=4745\", \"82\u00d763=5166\"],\n  [\"12\u00d795=1140\", \"78\u00d784=6552\"],\n  [\"79\u00d782=6478\", \"52\u00d792=4784\"],\n  [\"31\u00d783=2573\", \"39\u00d795=3705\"],\n  [\"47\u00d795=4465\", \"69\u00d760=4140\"],\n  [\"55\u00d744=2420\", \"22\u00d799=2178\"],\n  [\"82\u00d761=5002\", \"11\u00d725=275\"],\n  [\"67\u00d794=6298\", \"63\u00d772=4536\"],\n  [\"19\u00d784=1596\", \"21\u00d770=1470\"],\n  [\"50\u00d799=4950\", \"49\u00d767=3283\"],\n  [\"16\u00d771=1136\", \"92\u00d755=5060\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"33\u00d765=2145\"; new=\"83\u00d794=7802\"},\n    @{old=\"86\u00d725=2150\"; new=\"20\u00d780=1600\"},\n    @{old=\"92\u00d787=8004\"; new=\"65\u00d780=5200\"},\n    @{old=\"87\u00d793=8091\"; new=\"63\u00d741=2583\"},\n    @{old=\"97\u00d733=3201\"; new=\"34\u00d798=3332\"},\n    @{old=\"17\u00d794=1598\"; new=\"84\u00d772=6048\"},\n    @{old=\"55\u00d794=5170\"; new=\"12\u00d752=624\"},\n    @{old=\"59\u00d717=1003\"; new=\"45\u00d784=3780\"},\n    @{old=\"57\u00d749=2793\"; new=\"98\u00d790=8820\"},\n    @{old=\"70\u00d742=2940\"; new=\"29\u00d748=1392\"},\n    @{old=\"73\u00d753=3869\"; new=\"36\u00d760=2160\"},\n    @{old=\"92\u00d762=5704\"; new=\"97\u00d734=3298\"},\n    @{old=\"32\u00d713=416\"; new=\"25\u00d782=2050\"},\n    @{old=\"11\u00d770=770\"; new=\"65\u00d758=3770\"},\n    @{old=\"73\u00d765=4745\"; new=\"82\u00d763=5166\"},\n    @{old=\"12\u00d795=1140\"; new=\"78\u00d784=6552\"},\n    @{old=\"79\u00d782=6478\"; new=\"52\u00d792=4784\"},\n    @{old=\"31\u00d783=2573\"; new=\"39\u00d795=3705\"},\n    @{old=\"47\u00d795=4465\"; new=\"69\u00d760=4140\"},\n    @{old=\"55\u00d744=2420\"; new=\"22\u00d799=2178\"},\n    @{old=\"82\u00d761=5002\"; new=\"11\u00d725=275\"},\n    @{old=\"67\u00d794=6298\"; new=\"63\u00d772=4536\"},\n    @{old=\"19\u00d784=1596\"; new=\"21\u00d770=1470\"},\n    @{old=\"50\u00d799=4950\"; new=\"49\u00d767=3283\"},\n    @{old=\"16\u00d771=1136\"; new=\"92\u00d755=5060\"}\n)\n\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Execute(\n        $pair.old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        0,           # Wrap (wdFindStop)\n        $false,      # Format\n        $pair.new,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
